$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.469524
$ws.Range("H2").Value = 7.408571999999999
$ws.Range("I2").Value = 0.006775482240913427
$ws.Range("J2").Value = 0.006775482240913427
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 150.7437118387187
$ws.Range("R2").Value = 1356.693406548468
$ws.Range("S2").Value = 0.001384646670941439
$ws.Range("T2").Value = 0.001384646670941439

$ws.Range("G3").Value = 2.469524
$ws.Range("H3").Value = 7.408571999999999
$ws.Range("I3").Value = 0.006775482240913427
$ws.Range("J3").Value = 0.006775482240913427
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 262.546125334184
$ws.Range("R3").Value = 2362.915128007656
$ws.Range("S3").Value = 0.002411600550220613
$ws.Range("T3").Value = 0.002411600550220614

$ws.Range("G4").Value = 2.469524
$ws.Range("H4").Value = 7.408571999999999
$ws.Range("I4").Value = 0.006775482240913427
$ws.Range("J4").Value = 0.006775482240913427
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 324.3433539704906
$ws.Range("R4").Value = 2919.090185734416
$ws.Range("S4").Value = 0.002979235019751374
$ws.Range("T4").Value = 0.002979235019751375

$ws.Range("H5").Value = 988.862762
$ws.Range("I5").Value = 0.9043607975506752
$ws.Range("J5").Value = 0.9043607975506752
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 20120.59047856557
$ws.Range("R5").Value = 181085.3143070901
$ws.Range("S5").Value = 0.1848163899090481
$ws.Range("T5").Value = 0.1848163899090481

$ws.Range("H6").Value = 988.862762
$ws.Range("I6").Value = 0.9043607975506752
$ws.Range("J6").Value = 0.9043607975506752
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.321889559949188
$ws.Range("T6").Value = 0.3218895599491881

$ws.Range("H7").Value = 988.862762
$ws.Range("I7").Value = 0.9043607975506752
$ws.Range("J7").Value = 0.9043607975506752
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 43291.88740334885
$ws.Range("R7").Value = 389626.9866301397
$ws.Range("S7").Value = 0.397654847692439
$ws.Range("T7").Value = 0.3976548476924391

$ws.Range("G8").Value = 32.38899933333333
$ws.Range("H8").Value = 97.16699799999999
$ws.Range("I8").Value = 0.08886372020841134
$ws.Range("J8").Value = 0.08886372020841135
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 1977.076546835929
$ws.Range("R8").Value = 17793.68892152336
$ws.Range("S8").Value = 0.01816030947746387
$ws.Range("T8").Value = 0.01816030947746387

$ws.Range("G9").Value = 32.38899933333333
$ws.Range("H9").Value = 97.16699799999999
$ws.Range("I9").Value = 0.08886372020841134
$ws.Range("J9").Value = 0.08886372020841135
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 3443.419168397689
$ws.Range("R9").Value = 30990.7725155792
$ws.Range("S9").Value = 0.03162930532902768
$ws.Range("T9").Value = 0.03162930532902768

$ws.Range("G10").Value = 32.38899933333333
$ws.Range("H10").Value = 97.16699799999999
$ws.Range("I10").Value = 0.08886372020841134
$ws.Range("J10").Value = 0.08886372020841135
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 4253.919652338393
$ws.Range("R10").Value = 38285.27687104554
$ws.Range("S10").Value = 0.03907410540191979
$ws.Range("T10").Value = 0.0390741054019198
